$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list data (Price column D, Volume(1h) column E).
# D values are forced to Text (NumberFormat "@") before assignment so that
# numeric-looking strings (e.g. "0.9996", "10.07") are stored as text, matching
# the original inline-string cell type; ClearFormats() afterwards drops the
# temporary number-format style so the cell style index is left unchanged (0).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.929.99"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.68"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5055"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -3.37%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8937"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.99"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07441"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.05"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.231"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008476"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9995"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.962.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.015"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.089.51"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.432"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.90"
$ws.Range("D25").ClearFormats()

$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.063"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.16"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.664"
$ws.Range("D30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.672"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09234"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05082"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.993"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7461"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.149"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.282"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.531"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02000"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.084"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5350"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "117.99"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.498"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.530"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1467"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4651"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9988"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.561"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.84"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.09"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.73%  "
